# Bill of Materials: add a new "RF cable SMA Male to SMA Female 15cm" line
# right under the "HopeRF RFM95W module" row (new row 8), pushing every
# following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# 1. Insert a new row at 8 - Excel shifts rows 8:35 down to 9:36 and the
#    new row inherits row 7's formatting (style 4 in col B, style 1 in col C).
$ws.Rows("8:8").Insert()

# 2. Populate the new row.
$ws.Range("B8").Value = "RF cable SMA Male to SMA Female 15cm"
$ws.Range("C8").Value = 1

# 3. Hyperlink anchors are not shifted automatically by the row insert, so
#    drop all the (now stale) hyperlinks and re-create them against the
#    shifted cells, keeping their original target URLs, then add the new
#    hyperlink for the newly-inserted cable row.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B9"),  "https://www.st.com/en/evaluation-tools/nucleo-l476rg.html")
$ws.Hyperlinks.Add($ws.Range("B20"), "https://www.microchip.com/en-us/product/MCP1700")
$ws.Hyperlinks.Add($ws.Range("B7"),  "https://www.hoperf.com/modules/lora/RFM95.html")
$ws.Hyperlinks.Add($ws.Range("B19"), "https://www.pololu.com/product/1408")
$ws.Hyperlinks.Add($ws.Range("B12"), "https://www.infineon.com/cms/en/product/power/mosfet/n-channel/irlz44n/")
$ws.Hyperlinks.Add($ws.Range("B15"), "https://atlas-scientific.com/kits/ph-kit/")
$ws.Hyperlinks.Add($ws.Range("B16"), "https://atlas-scientific.com/kits/dissolved-oxygen-kit/")
$ws.Hyperlinks.Add($ws.Range("B14"), "https://atlas-scientific.com/kits/pt-1000-temperature-kit/")
$ws.Hyperlinks.Add($ws.Range("B8"),  "https://www.rfconnector.com/rf-cable/sma-male-to-sma-female-15cm")

# 4. Adding a hyperlink re-styles the cell with a freshly duplicated
#    "Hyperlink" xf; reassign the named style so it dedupes back onto the
#    workbook's existing Hyperlink style (same one the other description
#    cells already use).
$ws.Range("B7").Style = "Hyperlink"
$ws.Range("B8").Style = "Hyperlink"
$ws.Range("B9").Style = "Hyperlink"
$ws.Range("B12").Style = "Hyperlink"
$ws.Range("B14").Style = "Hyperlink"
$ws.Range("B15").Style = "Hyperlink"
$ws.Range("B16").Style = "Hyperlink"
$ws.Range("B19").Style = "Hyperlink"
$ws.Range("B20").Style = "Hyperlink"

# 5. Column B widens slightly to fit the new, longer description text.
$ws.Columns("B:B").ColumnWidth = 34.9

# 6. Restore the on-screen selection to where the editor left off.
$ws.Range("H6").Select()
